$d = $word.ActiveDocument

# Fix 1: remove the trailing space run after "${optionBlock}" (3rd paragraph).
# Target only the single stray space character (and not the paragraph mark or
# the preceding "}" run) so the surrounding runs stay untouched.
$p = $d.Paragraphs(3).Range
$full = $p.Text
$spaceStart = $p.Start + $full.Length - 2
$spaceEnd = $spaceStart + 1
$spaceRange = $d.Range($spaceStart, $spaceEnd)
$spaceRange.Delete()

Write-Host "done"
